# Notes.xlsx update (#LongDB):
#  - Log two more "Bảng Request" changes: the new Payment and RequestCode
#    columns, re-using the two still-blank template rows (5 & 6).
#  - Center the whole table body horizontally (and the "Nội dung" column
#    both horizontally + vertically).
#  - Widen the "Nội dung" column so the longer notes fit.
#  - Move the active-cell selection down to E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Fill rows 5 & 6 (previously blank placeholder rows) with the two new
#    change-log entries. Copy row 4's formatting first so the new rows pick
#    up the same borders / date number-format as the other data rows.
# ---------------------------------------------------------------------------
$ws.Range("B4:E4").Copy()
$ws.Range("B5:E6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B5").Value = 41717
$ws.Range("C5").Value = "Database"
$ws.Range("D5").Value = "Bảng Request"
$ws.Range("E5").Value = "Thêm cột Payment (bit-request đã được thanh toán chưa)"

$ws.Range("B6").Value = 41717
$ws.Range("C6").Value = "Database"
$ws.Range("D6").Value = "Bảng Request"
$ws.Range("E6").Value = "Thêm cột RequestCode (nvarchar-code cho staff check để thanh toán)"

# ---------------------------------------------------------------------------
# 2. Center-align the table body. Set alignment on one representative cell
#    per formatting group, then fan that format out with PasteSpecial so the
#    whole group shares a single style (instead of fragmenting cellXfs).
# ---------------------------------------------------------------------------

# B2:B6 - dated rows
$ws.Range("B2").HorizontalAlignment = -4108 # xlCenter
$ws.Range("B2").Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# B7:B29 - still-blank date-column cells
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").Copy()
$ws.Range("B7:B29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# B30 - bottom border row
$ws.Range("B30").HorizontalAlignment = -4108

# C2:D29 - Database / table-name columns (data + blank rows share one style)
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").Copy()
$ws.Range("C2:D29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# C30:D30 - bottom border row
$ws.Range("C30").HorizontalAlignment = -4108
$ws.Range("C30").Copy()
$ws.Range("C30:D30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# E2:E29 - Nội dung column (data + blank rows share one style), h + v center
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("E2").Copy()
$ws.Range("E2:E29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# E30 - bottom border row
$ws.Range("E30").HorizontalAlignment = -4108
$ws.Range("E30").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Widen column E to fit the longer text, and move the selection.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 83

$ws.Range("E9").Select()
